# Weekly fruit/vegetable price update.
# A new weekly record (row) is inserted into the "Espinaca" (spinach) price
# log, immediately after the existing row 506, pushing every subsequent
# record down by one row (old row 507 -> new row 508, ..., old row 539 ->
# new row 540).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 507; this shifts rows 507..539 down to
# 508..540 and extends the sheet dimension to A1:R540 automatically.
$ws.Rows.Item(507).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A507").Value = 9
$ws.Range("B507").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C507").Value = "Metropolitana"
$ws.Range("D507").Value = 44931
$ws.Range("D507").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E507").Value = 13
$ws.Range("F507").Value = 100112012
$ws.Range("G507").Value = "Espinaca"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 160
$ws.Range("K507").Value = 7000
$ws.Range("L507").Value = 9000
$ws.Range("M507").Value = 8000
$ws.Range("N507").Value = "$/cuna 10 kilos"
$ws.Range("O507").Value = "Provincia de Chacabuco"
$ws.Range("P507").Value = 800
$ws.Range("Q507").Value = 10
$ws.Range("R507").Value = "Hortaliza"
